# Scheduled data refresh: update crafting-profit computations in the
# Tonberry_Profits sheets (one tab per crafting class) with fresh market
# data. Only numeric result columns (H:N) change; labels/ids are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 40000340
$ws.Range("I11").Value = 40000340
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 40000340
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -40000200

# Row 19
$ws.Range("H19").Value = 2112.7273
$ws.Range("I19").Value = 572
$ws.Range("J19").Value = 2993.1428
$ws.Range("K19").Value = 572
$ws.Range("L19").Value = 2993.1428
$ws.Range("M19").Value = -397
$ws.Range("N19").Value = -3343.1428

# Row 43
$ws.Range("H43").Value = 1330.3334
$ws.Range("I43").Value = 999
$ws.Range("J43").Value = 1396.6
$ws.Range("K43").Value = 999
$ws.Range("L43").Value = 1396.6
$ws.Range("M43").Value = -930
$ws.Range("N43").Value = -1534.6

# Row 53
$ws.Range("H53").Value = 7889.2
$ws.Range("I53").Value = 10638.909
$ws.Range("J53").Value = 327.5
$ws.Range("K53").Value = 10638.909
$ws.Range("L53").Value = 327.5
$ws.Range("M53").Value = -10001.909
$ws.Range("N53").Value = -1601.5

# Row 62
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -4248

# Row 65
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -21240

# Row 92
$ws.Range("H92").Value = 1759433.2
$ws.Range("I92").Value = 2052506
$ws.Range("J92").Value = 997
$ws.Range("K92").Value = 2052506
$ws.Range("L92").Value = 997
$ws.Range("M92").Value = -2051258
$ws.Range("N92").Value = -3493

# Row 96
$ws.Range("H96").Value = 2342
$ws.Range("I96").Value = 2342
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 7026
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -5653

# Row 97
$ws.Range("H97").Value = 1192.6666
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1192.6666
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3577.9998
$ws.Range("N97").Value = -4569.9998

# Row 98
$ws.Range("H98").Value = 2202.25
$ws.Range("I98").Value = 2538.5386
$ws.Range("J98").Value = 745
$ws.Range("K98").Value = 2538.5386
$ws.Range("L98").Value = 745
$ws.Range("M98").Value = -1040.5386
$ws.Range("N98").Value = -3741

# Row 100
$ws.Range("H100").Value = 2821
$ws.Range("I100").Value = 2821
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2821
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2280

# Row 112
$ws.Range("H112").Value = 6999.5
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 7842.2856
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 23526.8568
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -25742.8568

# Row 122
$ws.Range("H122").Value = 2202.25
$ws.Range("I122").Value = 2538.5386
$ws.Range("J122").Value = 745
$ws.Range("K122").Value = 7615.6158
$ws.Range("L122").Value = 2235
$ws.Range("M122").Value = -5165.6158
$ws.Range("N122").Value = -7135

# Row 132
$ws.Range("H132").Value = 751.6491
$ws.Range("I132").Value = 697.6326
$ws.Range("J132").Value = 1082.5
$ws.Range("K132").Value = 2092.8978
$ws.Range("L132").Value = 3247.5
$ws.Range("M132").Value = 437.1021999999998
$ws.Range("N132").Value = -8307.5

# Row 137
$ws.Range("H137").Value = 1648.2727
$ws.Range("I137").Value = 1327.5385
$ws.Range("J137").Value = 2111.5557
$ws.Range("K137").Value = 3982.6155
$ws.Range("L137").Value = 6334.6671
$ws.Range("M137").Value = -1432.6155
$ws.Range("N137").Value = -11434.6671

# Row 138
$ws.Range("H138").Value = 1744.3485
$ws.Range("I138").Value = 1241
$ws.Range("J138").Value = 2163.8057
$ws.Range("K138").Value = 3723
$ws.Range("L138").Value = 6491.4171
$ws.Range("M138").Value = 1417
$ws.Range("N138").Value = -16771.4171

# Row 141
$ws.Range("H141").Value = 4564.467
$ws.Range("I141").Value = 3606.4443
$ws.Range("J141").Value = 6001.5
$ws.Range("K141").Value = 10819.3329
$ws.Range("L141").Value = 18004.5
$ws.Range("M141").Value = -5639.332900000001
$ws.Range("N141").Value = -28364.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3010.6934
$ws.Range("I32").Value = 1927.197
$ws.Range("J32").Value = 10956.333
$ws.Range("K32").Value = 1927.197
$ws.Range("L32").Value = 10956.333
$ws.Range("M32").Value = -1640.197
$ws.Range("N32").Value = -11530.333

# Row 45
$ws.Range("H45").Value = 1327.96
$ws.Range("I45").Value = 1013
$ws.Range("J45").Value = 1887.8889
$ws.Range("K45").Value = 1013
$ws.Range("L45").Value = 1887.8889
$ws.Range("M45").Value = -636
$ws.Range("N45").Value = -2641.8889

# Row 53
$ws.Range("H53").Value = 19500
$ws.Range("I53").Value = 19000
$ws.Range("J53").Value = 20000
$ws.Range("K53").Value = 19000
$ws.Range("L53").Value = 20000
$ws.Range("M53").Value = -18318
$ws.Range("N53").Value = -21364

# Row 61
$ws.Range("H61").Value = 3997.7307
$ws.Range("I61").Value = 2918
$ws.Range("J61").Value = 4923.2144
$ws.Range("K61").Value = 2918
$ws.Range("L61").Value = 4923.2144
$ws.Range("M61").Value = -2706
$ws.Range("N61").Value = -5347.2144

# Row 132
$ws.Range("H132").Value = 1571
$ws.Range("I132").Value = 906.0345
$ws.Range("J132").Value = 3499.4
$ws.Range("K132").Value = 2718.1035
$ws.Range("L132").Value = 10498.2
$ws.Range("M132").Value = -188.1035000000002
$ws.Range("N132").Value = -15558.2

# Row 136
$ws.Range("H136").Value = 3997.7307
$ws.Range("I136").Value = 2918
$ws.Range("J136").Value = 4923.2144
$ws.Range("K136").Value = 8754
$ws.Range("L136").Value = 14769.6432
$ws.Range("M136").Value = -6204
$ws.Range("N136").Value = -19869.6432

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1491.909
$ws.Range("I99").Value = 1401.5454
$ws.Range("J99").Value = 1672.6364
$ws.Range("K99").Value = 1401.5454
$ws.Range("L99").Value = 1672.6364
$ws.Range("M99").Value = 96.45460000000003
$ws.Range("N99").Value = -4668.6364

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2353.6191
$ws.Range("I31").Value = 2057.4546
$ws.Range("J31").Value = 2679.4
$ws.Range("K31").Value = 2057.4546
$ws.Range("L31").Value = 2679.4
$ws.Range("M31").Value = -1762.4546
$ws.Range("N31").Value = -3269.4

# Row 34
$ws.Range("H34").Value = 2353.6191
$ws.Range("I34").Value = 2057.4546
$ws.Range("J34").Value = 2679.4
$ws.Range("K34").Value = 2057.4546
$ws.Range("L34").Value = 2679.4
$ws.Range("M34").Value = -1855.4546
$ws.Range("N34").Value = -3083.4

# Row 134
$ws.Range("H134").Value = 1049.04
$ws.Range("I134").Value = 1022.5714
$ws.Range("J134").Value = 1188
$ws.Range("K134").Value = 3067.7142
$ws.Range("L134").Value = 3564
$ws.Range("M134").Value = -532.7142000000003
$ws.Range("N134").Value = -8634

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 33.75
$ws.Range("I8").Value = 33.75
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 101.25
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 37.75

# Row 116
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = ""

# Row 131
$ws.Range("H131").Value = 801.9299999999999
$ws.Range("I131").Value = 576.1667
$ws.Range("J131").Value = 816.34045
$ws.Range("K131").Value = 1728.5001
$ws.Range("L131").Value = 2449.02135
$ws.Range("M131").Value = 3311.4999
$ws.Range("N131").Value = -12529.02135

# Row 136
$ws.Range("H136").Value = 1734
$ws.Range("I136").Value = 1734
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5202
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -102

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 756.8125
$ws.Range("I93").Value = 761.3570999999999
$ws.Range("J93").Value = 725
$ws.Range("K93").Value = 761.3570999999999
$ws.Range("L93").Value = 725
$ws.Range("M93").Value = 486.6429000000001
$ws.Range("N93").Value = -3221

# Row 132
$ws.Range("H132").Value = 2450.25
$ws.Range("I132").Value = 2043.826
$ws.Range("J132").Value = 3169.3076
$ws.Range("K132").Value = 6131.478
$ws.Range("L132").Value = 9507.9228
$ws.Range("M132").Value = -3601.478
$ws.Range("N132").Value = -14567.9228

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2374.75
$ws.Range("I81").Value = 2249.5
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 4499
$ws.Range("L81").Value = 5000
$ws.Range("M81").Value = -3438
$ws.Range("N81").Value = -7122

# Row 84
$ws.Range("H84").Value = 2374.75
$ws.Range("I84").Value = 2249.5
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 22495
$ws.Range("L84").Value = 25000
$ws.Range("M84").Value = -17191
$ws.Range("N84").Value = -35608
